# documentation %26 message class
#
# The {R-STRING} marker cell (C5) becomes a generic {R-TEXT} marker.
# Setting a brand-new value on C5 naturally retires the old, now-unused
# "{R-STRING}" shared-string entry and appends "{R-TEXT}" as a new one,
# which is exactly what the OOXML diff shows happening to sharedStrings.xml
# (and, as a side effect, renumbers the other shared-string indices used
# by B6/C6/B8 — their text stays the same).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "{R-TEXT}"

# Reset the lingering UI selection (previously parked on D1) back to the
# sheet's home cell.
[void]$ws.Range("A1").Select()
